$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 56 (Steven Mustermann, 10a): the automatic logout time
#     got corrected / the student's departure is now flagged as "Abmeldung Automatisch"
#     at 16:00 instead of 18:00, shifting the (already negative/placeholder) duration.
$ws.Range("E56").Value = "08.06.2024 16:00"
$ws.Range("F56").Value = -329
$ws.Range("G56").Value = "Abmeldung Automatisch"

# --- Append new Anmeldehistorie rows for the 15.06.2024 20:54 / 18:00 event,
#     one row per affected student (Eli Enders 4a, Detlef Soost 1a, Tim Test 11c).
$ws.Range("A58").Value = "Eli"
$ws.Range("B58").Value = "Enders"
$ws.Range("C58").Value = "4a"
$ws.Range("D58").Value = "15.06.2024 20:54"
$ws.Range("E58").Value = "15.06.2024 18:00"
$ws.Range("F58").Value = -174

$ws.Range("A59").Value = "Detlef"
$ws.Range("B59").Value = "Soost"
$ws.Range("C59").Value = "1a"
$ws.Range("D59").Value = "15.06.2024 20:54"
$ws.Range("E59").Value = "15.06.2024 18:00"
$ws.Range("F59").Value = -174

$ws.Range("A60").Value = "Tim"
$ws.Range("B60").Value = "Test"
$ws.Range("C60").Value = "11c"
$ws.Range("D60").Value = "15.06.2024 20:54"
$ws.Range("E60").Value = "15.06.2024 18:00"
$ws.Range("F60").Value = -174
